$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.121.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.68%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.647.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.60%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +1.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'216.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.41%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.508"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.58%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.72%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.0638"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.25%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.254"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.84%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.76%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0798"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.38%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.712.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +5.35%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.01%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.542"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.42%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'63.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.78%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₃0761"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.43%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'25.975.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.08%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +0.90%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'194.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.67%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.75%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'9.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.66%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.69%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +1.35%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'144.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.17%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.34%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.95%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.12%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'15.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.02%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.18%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0489"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.82%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.66%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.01%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.56%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'0.901"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.13%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.131.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.57%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.67%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -2.56%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.94%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.799"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.06%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'98.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.71%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -2.93%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -1.53%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'56.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.53%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.48%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0523"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.75%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'Mantle"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.419"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.17%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'7.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.57%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.71%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -2.28%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +1.29%  "
$ws.Range("E51").Style = "Normal"
